# Apply edits to the "Manchester Utd_stats" workbook:
#  1. Bump the "age-days" value in column E (rows 4-44) by one day on every
#     per-player stats sheet (all sheets except "Matches").
#  2. Fix the "Playing Time" header block on the StandardStats and
#     PlayingTime sheets: it should occupy G1:I1 (with F1 holding the
#     "Unnamed: 4_level_0" placeholder used on every other stats sheet)
#     instead of F1:I1.
#  3. Rename the stats sheet tabs to their human readable names.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Increment the day-of-year portion of the "age-days" strings found in
#    column E of every stats sheet (everything except "Matches").
# ---------------------------------------------------------------------
$statsSheetNames = @(
    "StandardStats",
    "ShootingStats",
    "PassingStats",
    "PassTypes",
    "GoalShotCreation",
    "DefensiveActions",
    "Possession",
    "PlayingTime",
    "MiscStats"
)

foreach ($sheetName in $statsSheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    for ($row = 4; $row -le 44; $row++) {
        $cell = $ws.Range("E" + $row)
        $val = $cell.Value2
        if ($val -ne $null -and $val.Contains("-")) {
            $parts = $val.Split("-")
            if ($parts.Length -eq 2 -and $parts[1].Length -eq 3) {
                $year = $parts[0]
                $newDay = [int]$parts[1] + 1
                $newDayStr = ([string]$newDay).PadLeft(3, '0')
                $cell.Value = $year + "-" + $newDayStr
            }
        }
    }
}

# ---------------------------------------------------------------------
# 2. Fix the "Playing Time" merged header on StandardStats & PlayingTime.
#    Currently: F1:I1 merged, F1 = "Playing Time".
#    Target:    G1:I1 merged, F1 = "Unnamed: 4_level_0", G1 = "Playing Time".
# ---------------------------------------------------------------------
$headerSheetNames = @("StandardStats", "PlayingTime")

foreach ($sheetName in $headerSheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    # Clear borders before unmerging so the engine does not fabricate
    # extra "split border" style variants for the soon-to-be-discarded
    # layout.
    $ws.Range("F1:I1").Borders.LineStyle = 0
    $ws.Range("V1:X1").Borders.LineStyle = 0

    $ws.Range("F1:I1").UnMerge()
    $ws.Range("F1").Value = "Unnamed: 4_level_0"
    $ws.Range("G1").Value = "Playing Time"

    # Temporarily unmerge the next merged block so the new merge gets
    # inserted in the same slot the old one occupied (keeps mergeCells
    # ordering stable), then restore it.
    $ws.Range("V1:X1").UnMerge()
    $ws.Range("G1:I1").Merge()
    $ws.Range("V1:X1").Merge()

    # Restore the original header formatting (font, alignment, thin box
    # border) by copying it from a still-pristine header cell (A1 uses
    # the same style as the rest of row 1).
    $ws.Range("A1").Copy()
    $ws.Range("F1:I1").PasteSpecial(-4122)
    $ws.Range("V1:X1").PasteSpecial(-4122)
    $excel.CutCopyMode = 0
}

# ---------------------------------------------------------------------
# 3. Rename the stats sheet tabs.
# ---------------------------------------------------------------------
$renameMap = @{
    "StandardStats"    = "Standard Stats"
    "ShootingStats"    = "Shooting Stats"
    "PassingStats"     = "Passing Stats"
    "PassTypes"        = "Pass Types"
    "GoalShotCreation" = "Goal & Shot Creation"
    "DefensiveActions" = "Defensive Actions"
    "PlayingTime"      = "Playing Time"
    "MiscStats"        = "Miscellaneous Stats"
}

foreach ($oldName in $renameMap.Keys) {
    $ws = $wb.Worksheets.Item($oldName)
    $ws.Name = $renameMap[$oldName]
}
